# Apply the latest scraped crypto price/volume snapshot to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.660.11'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '''3.448.99'
$ws.Range("E3").Value = '  -3.20%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''593.35'
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("D6").Value = '''136.48'
$ws.Range("E6").Value = '  -7.27%  '
$ws.Range("D7").Value = '''3.447.73'
$ws.Range("E7").Value = '  -3.17%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '''0.499'
$ws.Range("E9").Value = '  +2.13%  '
$ws.Range("E10").Value = '  -5.70%  '
$ws.Range("E11").Value = '  -8.53%  '
$ws.Range("D12").Value = '''0.378'
$ws.Range("E12").Value = '  -7.82%  '
$ws.Range("D13").Value = '''4.032.09'
$ws.Range("E13").Value = '  -3.26%  '
$ws.Range("E14").Value = '  -10.26%  '
$ws.Range("D15").Value = '''26.64'
$ws.Range("E15").Value = '  -9.12%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '''3.441.08'
$ws.Range("E16").Value = '  -3.18%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '''65.607.33'
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '''0.115'
$ws.Range("E18").Value = '  -2.14%  '
$ws.Range("D19").Value = '''9.88'
$ws.Range("E19").Value = '  -10.78%  '
$ws.Range("D20").Value = '''5.83'
$ws.Range("E20").Value = '  -6.99%  '
$ws.Range("D21").Value = '''13.75'
$ws.Range("E21").Value = '  -7.14%  '
$ws.Range("D22").Value = '''394.09'
$ws.Range("E22").Value = '  -6.66%  '
$ws.Range("D23").Value = '''0.552'
$ws.Range("E23").Value = '  -9.01%  '
$ws.Range("D24").Value = '''73.55'
$ws.Range("E24").Value = '  -5.64%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = '''3.590.18'
$ws.Range("E26").Value = '  -3.07%  '
$ws.Range("D27").Value = '''0.0000107'
$ws.Range("E27").Value = '  -9.82%  '
$ws.Range("D28").Value = '''0.997'
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''2.26'
$ws.Range("E29").Value = '  -9.17%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '''7.20'
$ws.Range("E30").Value = '  -9.33%  '
$ws.Range("D31").Value = '''8.22'
$ws.Range("E31").Value = '  -11.22%  '
$ws.Range("D32").Value = '''3.454.46'
$ws.Range("E32").Value = '  -2.98%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("E34").Value = '  -5.97%  '
$ws.Range("D35").Value = '''23.02'
$ws.Range("E35").Value = '  -7.04%  '
$ws.Range("D36").Value = '''171.96'
$ws.Range("E36").Value = '  -1.58%  '
$ws.Range("D37").Value = '''6.97'
$ws.Range("E37").Value = '  -9.14%  '
$ws.Range("D38").Value = '''1.20'
$ws.Range("E38").Value = '  -10.24%  '
$ws.Range("E39").Value = '  -6.83%  '
$ws.Range("D40").Value = '''4.83'
$ws.Range("E40").Value = '  -10.30%  '
$ws.Range("E41").Value = '  -7.42%  '
$ws.Range("E42").Value = '  -4.75%  '
$ws.Range("D43").Value = '''43.62'
$ws.Range("E43").Value = '  -4.78%  '
$ws.Range("D44").Value = '''1.00'
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").Value = '''4.42'
$ws.Range("E45").Value = '  -14.00%  '
$ws.Range("D46").Value = '''1.63'
$ws.Range("E46").Value = '  -11.50%  '
$ws.Range("E47").Value = '  -1.52%  '
$ws.Range("D48").Value = '''22.42'
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("D49").Value = '''6.55'
$ws.Range("E49").Value = '  -7.80%  '
$ws.Range("E50").Value = '  -14.79%  '
$ws.Range("D51").Value = '''2.202.62'
$ws.Range("E51").Value = '  -7.67%  '
